$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Inhbb","Acvr2a","ECs",2,0.6666666666666666,3.345805333333333,10.037416,0.3489465220682754,0.3489465220682754,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,26.56704047567911,239.103364281112,0.06740421970872093,0.06740421970872093),
    @("ECs","Inhbb","Acvr2a","FAPs",2,0.6666666666666666,3.345805333333333,10.037416,0.3489465220682754,0.3489465220682754,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,90.18232616171909,811.6409354554719,0.2288049108071629,0.2288049108071629),
    @("ECs","Inhbb","Acvr2a","sCs",2,0.6666666666666666,3.345805333333333,10.037416,0.3489465220682754,0.3489465220682754,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,20.78618255665155,187.075643009864,0.05273739155239154,0.05273739155239155),
    @("FAPs","Inhbb","Acvr2a","ECs",3,1,5.480061666666667,16.440185,0.5715360783999618,0.5715360783999618,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,43.51389444481056,391.625050003295,0.1104007088868309,0.1104007088868309),
    @("FAPs","Inhbb","Acvr2a","FAPs",3,1,5.480061666666667,16.440185,0.5715360783999618,0.5715360783999618,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,147.7087455405855,1329.37870986527,0.3747573142906758,0.3747573142906758),
    @("FAPs","Inhbb","Acvr2a","sCs",3,1,5.480061666666667,16.440185,0.5715360783999618,0.5715360783999618,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,34.04548408426277,306.409356758365,0.08637805522245509,0.0863780552224551),
    @("sCs","Inhbb","Acvr2a","ECs",3,1,0.7624369999999999,2.287311,0.07951739953176286,0.07951739953176286,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,6.054056533819667,54.486508804377,0.01535997045316984,0.01535997045316984),
    @("sCs","Inhbb","Acvr2a","FAPs",3,1,0.7624369999999999,2.287311,0.07951739953176286,0.07951739953176286,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,20.55061049928466,184.955494493562,0.05213971298422249,0.05213971298422249),
    @("sCs","Inhbb","Acvr2a","sCs",3,1,0.7624369999999999,2.287311,0.07951739953176286,0.07951739953176286,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,4.736723476424332,42.630511287819,0.01201771609437053,0.01201771609437053)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
